$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = @(0.675330155194221903869333800685, 0.312790395851139102489923971007, 0.805489636583999191188354416227, 8.660232485948974101575004169717, 10.453842673578330746408937557135)
    3 = @(0.000000207022599729711487270833, 0.002777888934908601115125748038, 0.805489636583999191188354416227, 0.496779210170732010709571113694, 1.305046942712239976813748398854)
    4 = @(0.000858366962651846421294976608, 0.002777888934908601115125748038, 0.805489636583999191188354416227, 8.660232485948974101575004169717, 9.469358378430534628478199010715)
    5 = @(0.000020749860322855081471707037, 0.000070973895028636491133511299, 0.157525292976961495128307433333, 8.660232485948974101575004169717, 8.817849502681287532368514803238)
    6 = @(0.003994804209775715264640894020, 0.312790395851139102489923971007, 0.157525292976961495128307433333, 0.496779210170732010709571113694, 0.971089703208608323592443412053)
}

foreach ($row in $values.Keys) {
    $vals = $values[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("G$row").Value = $vals[4]
}
